# NetworkDictionary.xlsx: update templates to fit the new data model.
#
# The "CollectionEvents" sheet previously described an age range with two
# separate columns ("ageMin", "ageMax"). The data model now instead uses a
# single "ageGroups" column, so we drop the "ageMin" column and rename the
# former "ageMax" header to "ageGroups" (which, after the shift, lands in
# the same column the old "ageMax" occupied).

$wb = $excel.ActiveWorkbook

# Remember/restore the active sheet so selecting ranges on other sheets
# below doesn't change which tab is marked active.
$originalActive = $wb.ActiveSheet

$wsCollectionEvents = $wb.Worksheets.Item("CollectionEvents")

# Drop the "ageMin" column (column D). This shifts the old "ageMax" (E) and
# "subcohorts" (F) columns left by one, onto D and E respectively - and
# carries the "subcohorts" column's width/formatting onto its new position.
$wsCollectionEvents.Columns.Item(4).Delete()

# The cell that used to read "ageMax" now holds the new "ageGroups" header.
$wsCollectionEvents.Range("D1").Value = "ageGroups"

# Approximate the original column's best-fit width for the new header text.
$wsCollectionEvents.Columns.Item(4).ColumnWidth = 9

# Restore the view state (selection) seen in the authored workbook.
$wsTargetTables = $wb.Worksheets.Item("TargetTables")
$wsTargetTables.Range("C12:C13").Select()

$wsCollectionEvents.Range("F3").Select()

$originalActive.Activate()
